$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.518.50"
$ws.Range("E2").Value = "  +2.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.812.28"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "681.41"
$ws.Range("E5").Value = "  +8.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.66"
$ws.Range("E6").Value = "  +3.97%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.811.54"
$ws.Range("E7").Value = "  +1.02%  "
$ws.Range("E8").Value = "  -0.18%  "
$ws.Range("E9").Value = "  +0.84%  "
$ws.Range("E10").Value = "  +1.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.29"
$ws.Range("E11").Value = "  +7.05%  "
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000246"
$ws.Range("E13").Value = "  +0.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.99"
$ws.Range("E14").Value = "  +2.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.455.87"
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.812.21"
$ws.Range("E16").Value = "  +1.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70.582.08"
$ws.Range("E17").Value = "  +2.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.71"
$ws.Range("E18").Value = "  +0.77%  "
$ws.Range("E19").Value = "  +2.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.115"
$ws.Range("E20").Value = "  +0.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.28"
$ws.Range("E21").Value = "  +18.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "477.83"
$ws.Range("E22").Value = "  +2.50%  "
$ws.Range("E23").Value = "  +1.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.49"
$ws.Range("E24").Value = "  +0.61%  "
$ws.Range("E25").Value = "  -1.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.26"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.34"
$ws.Range("E27").Value = "  +3.21%  "
$ws.Range("E28").Value = "  -1.23%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.964.59"
$ws.Range("E30").Value = "  +0.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.92"
$ws.Range("E31").Value = "  +9.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.30"
$ws.Range("E32").Value = "  +2.91%  "
$ws.Range("E33").Value = "  +4.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.62"
$ws.Range("E34").Value = "  +3.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.180"
$ws.Range("E35").Value = "  +4.65%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "Binance-PegBSC-USD"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.44%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.14"
$ws.Range("E37").Value = "  +2.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.769.83"
$ws.Range("E38").Value = "  +1.16%  "
$ws.Range("E39").Value = "  +1.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.40"
$ws.Range("E40").Value = "  +3.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.96"
$ws.Range("E41").Value = "  +2.48%  "
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.15"
$ws.Range("E44").Value = "  +13.60%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "46.08"
$ws.Range("E46").Value = "  +7.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "159.99"
$ws.Range("E47").Value = "  +2.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.47"
$ws.Range("E48").Value = "  +8.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "48.08"
$ws.Range("E49").Value = "  +3.14%  "
$ws.Range("E50").Value = "  +9.38%  "
$ws.Range("E51").Value = "  +1.71%  "
